$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 <-> row 4, and row 3 <-> row 5
# across columns D, M, N, O, P, R, S (other columns are identical between
# the swapped rows so they are left untouched).

$row2 = @{ D = 44362; M = 100; N = 19000; O = 20000; P = 19500; R = "Provincia de Curicó";  S = 1083 }
$row3 = @{ D = 45084; M = 100; N = 17000; O = 18000; P = 17500; R = "Región de O'Higgins";  S = 972  }
$row4 = @{ D = 44320; M = 50;  N = 18000; O = 20000; P = 18800; R = "Provincia de Limarí";  S = 1044 }
$row5 = @{ D = 44719; M = 50;  N = 20000; O = 21000; P = 20400; R = "Provincia de Limarí";  S = 1133 }

$updates = @{ 2 = $row2; 3 = $row3; 4 = $row4; 5 = $row5 }

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
